# Update "paises.xlsx" (Pais worksheet) with refreshed COVID-19 country
# statistics and the corresponding re-sort of rows whose total-case count
# changed enough to overtake a neighboring row.
#
# Column layout (row 3 = header):
#   A Pais | B Casos totales | C Nuevos casos | D Casos activos
#   E Recuperados | F Casos criticos | G Muertes hoy | H Muertes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $pais, $totales, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value = $pais
    $ws.Cells.Item($row, 2).Value = $totales
    $ws.Cells.Item($row, 3).Value = $nuevos
    $ws.Cells.Item($row, 4).Value = $activos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $criticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# --- Timestamp banner (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 15:52"

# --- Alemania overtakes Francia (rows 7-8) ---------------------------------
Set-Row 7 "Alemania" 109178 1515 36081 71031 4895 50 2066
Set-Row 8 "Francia" 109069 0 19337 79404 7131 0 10328

# --- Rumania data refresh (row 32, no re-order) -----------------------------
Set-Row 32 "Rumania" 4761 344 528 4018 162 18 215

# --- Pakistan overtakes Malasia (rows 34-35) --------------------------------
Set-Row 34 "Pakistan" 4183 148 467 3658 25 1 58
Set-Row 35 "Malasia" 4119 156 1487 2567 76 2 65

# --- Emiratos Arabes Unidos jumps above Finlandia & Tailandia (rows 44-46) --
Set-Row 44 "Emiratos Arabes Unidos" 2659 300 239 2408 1 0 12
Set-Row 45 "Finlandia" 2487 179 300 2147 82 6 40
Set-Row 46 "Tailandia" 2369 111 888 1451 61 3 30

# --- Islandia data refresh (row 56, no re-order) ----------------------------
Set-Row 56 "Islandia" 1616 30 633 977 13 0 6

# --- Senegal data refresh (row 109, no re-order) ----------------------------
Set-Row 109 "Senegal" 244 7 113 129 1 0 2

# --- Sri Lanka data refresh (row 113, no re-order) --------------------------
Set-Row 113 "Sri Lanka" 189 4 44 138 5 1 7

# --- Jamaica data refresh (row 139, no re-order) ----------------------------
Set-Row 139 "Jamaica" 63 0 10 50 0 0 3
